# "modify year option to all forms"
# Insert a new column before column B and add the "Year of Competition"
# header, shifting all the existing header columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B (existing B..AO shift right to C..AP),
# carrying formatting from the left (column A).
$ws.Columns("B").Insert() | Out-Null

# Match column A's width for the freshly inserted column.
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Populate the new header cell.
$ws.Range("B1").Value = "比賽年份 Year of Competition"

# Leave selection where Excel would naturally land after the insert.
$ws.Range("C2").Select() | Out-Null
